$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price / 1h-volume-change data to match the latest feed.
#
# Some "Price" cells contain plain-looking decimal strings (e.g. "23.26").
# Assigning such a string straight to .Value lets Excel auto-convert it into a
# real number, which would change the cell's stored type away from the
# original text. To keep those specific cells as literal text we briefly force
# the cell to a Text number format before the write and restore the original
# "Normal" style immediately after, so no extra formatting is left behind.
function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '27.714.59'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '1.645.53'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.24%  '
Set-TextValue $ws 'D5' '213.26'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  +0.22%  '
Set-TextValue $ws 'D8' '23.26'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('E10').Value = '  +0.56%  '
Set-TextValue $ws 'D11' '0.0894'
$ws.Range('E11').Value = '  +0.28%  '
$ws.Range('D12').Value = '1.878.70'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').Value = '1.647.40'
$ws.Range('E13').Value = '  +0.35%  '
Set-TextValue $ws 'D15' '0.557'
$ws.Range('E15').Value = '  -0.44%  '
Set-TextValue $ws 'D16' '64.66'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = '27.704.32'
$ws.Range('E17').Value = '  +1.12%  '
Set-TextValue $ws 'D18' '231.44'
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('E22').Value = '  -0.73%  '
Set-TextValue $ws 'D23' '10.08'
$ws.Range('E23').Value = '  +10.19%  '
Set-TextValue $ws 'D24' '1.96'
$ws.Range('E24').Value = '  -3.79%  '
Set-TextValue $ws 'D25' '150.32'
$ws.Range('E25').Value = '  +1.82%  '
Set-TextValue $ws 'D26' '6.91'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  -2.63%  '
$ws.Range('E28').Value = '  +0.24%  '
Set-TextValue $ws 'D29' '15.62'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('E31').Value = '  -0.22%  '
Set-TextValue $ws 'D32' '3.30'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').Value = '1.441.31'
$ws.Range('E33').Value = '  +3.30%  '
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('E37').Value = '  +1.26%  '
Set-TextValue $ws 'D38' '0.877'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('E40').Value = '  +12.78%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D43' '67.09'
$ws.Range('E43').Value = '  +4.43%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D44' '5.57'
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D45' '2.25'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.788.62'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D47' '1.74'
$ws.Range('E47').Value = '  +5.64%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0108'
$ws.Range('E48').Value = '  +2.94%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D49' '85.44'
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D50' '0.0987'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D51' '7.73'
$ws.Range('E51').Value = '  +1.58%  '
